# Update "想去人数" (want-to-go count) values in column F for the 展览
# (Exhibition) sheet and the 全部类型 (All types) combined sheet, matching
# the latest scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 821
$ws1.Range("F4").Value = 1403
$ws1.Range("F5").Value = 837
$ws1.Range("F6").Value = 476
$ws1.Range("F7").Value = 627
$ws1.Range("F9").Value = 7
$ws1.Range("F10").Value = 48
$ws1.Range("F11").Value = 210
$ws1.Range("F12").Value = 115
$ws1.Range("F13").Value = 1603
$ws1.Range("F14").Value = 202
$ws1.Range("F15").Value = 32
$ws1.Range("F17").Value = 71
$ws1.Range("F18").Value = 394
$ws1.Range("F19").Value = 109
$ws1.Range("F21").Value = 30
$ws1.Range("F22").Value = 216
$ws1.Range("F25").Value = 1452
$ws1.Range("F26").Value = 172

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 821
$ws4.Range("F5").Value = 1403
$ws4.Range("F6").Value = 837
$ws4.Range("F9").Value = 476
$ws4.Range("F10").Value = 627
$ws4.Range("F13").Value = 7
$ws4.Range("F14").Value = 48
$ws4.Range("F15").Value = 210
$ws4.Range("F16").Value = 115
$ws4.Range("F17").Value = 1603
$ws4.Range("F19").Value = 202
$ws4.Range("F20").Value = 32
$ws4.Range("F22").Value = 71
$ws4.Range("F23").Value = 394
$ws4.Range("F25").Value = 109
$ws4.Range("F33").Value = 30
$ws4.Range("F34").Value = 216
$ws4.Range("F37").Value = 1452
$ws4.Range("F38").Value = 172
